$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 158; this shifts the existing rows 158-235
# down to 159-236 (and the sheet dimension grows from R235 to R236).
$ws.Rows(158).Insert()

# Populate the newly inserted row 158 with the new "Acelga" weekly record.
$ws.Range("A158").Value = 4
$ws.Range("B158").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C158").Value = "Los Lagos"
$ws.Range("D158").Value = 44873
$ws.Range("E158").Value = 10
$ws.Range("F158").Value = 100112009
$ws.Range("G158").Value = "Acelga"
$ws.Range("H158").Value = "Sin especificar"
$ws.Range("I158").Value = "Primera"
$ws.Range("J158").Value = 200
$ws.Range("K158").Value = 3000
$ws.Range("L158").Value = 3000
$ws.Range("M158").Value = 3000
$ws.Range("N158").Value = '$/docena de atados (4 kilos)'
$ws.Range("O158").Value = "Región del Maule"
$ws.Range("P158").Value = 750
$ws.Range("Q158").Value = 4
$ws.Range("R158").Value = "Hortaliza"
